# Add the new "Prompt2" value for the Finance row (row 2, column D)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "What are the sales record for past year quarter wise"

# Match column D's resulting best-fit width as closely as this runtime allows
$ws.Columns.Item(4).ColumnWidth = 41.7

# Leave the active selection on D3, matching the saved view state
$ws.Range("D3").Select() | Out-Null
